$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the room for each session from "Uppsala" to "Online"
$ws.Range("B2").Value = "Online"
$ws.Range("B8").Value = "Online"
$ws.Range("B15").Value = "Online"
$ws.Range("B19").Value = "Online"
$ws.Range("B24").Value = "Online"

# Update the active selection as left by the author
$ws.Range("B29").Select()
